$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43
$ws.Range("A43").Value = "2025-12-03T09:52:08.671Z"
$ws.Range("B43").Value = "test@test.com"
$ws.Range("C43").Value = "https://www.mim.gov.sa/ar"
$ws.Range("D43").Value = "No"
$ws.Range("E43").Value = "٣‏/١٢‏/٢٠٢٥"
$ws.Range("F43").Value = "١٢:٥٢:٠٨ م"

# Row 44
$ws.Range("A44").Value = "2025-12-03T09:58:02.766Z"
$ws.Range("B44").Value = "test@test.com"
$ws.Range("C44").Value = "https://www.mim.gov.sa/ar"
$ws.Range("D44").Value = "No"
$ws.Range("E44").Value = "٣‏/١٢‏/٢٠٢٥"
$ws.Range("F44").Value = "١٢:٥٨:٠٢ م"

# Row 45
$ws.Range("A45").Value = "2025-12-03T09:58:37.673Z"
$ws.Range("B45").Value = "test@test.com"
$ws.Range("C45").Value = "https://www.mim.gov.sa/ar"
$ws.Range("D45").Value = "No"
$ws.Range("E45").Value = "٣‏/١٢‏/٢٠٢٥"
$ws.Range("F45").Value = "١٢:٥٨:٣٧ م"

# Row 46
$ws.Range("A46").Value = "2025-12-03T10:00:57.299Z"
$ws.Range("B46").Value = "test@test.com"
$ws.Range("C46").Value = "https://www.mim.gov.sa/ar"
$ws.Range("D46").Value = "No"
$ws.Range("E46").Value = "٣‏/١٢‏/٢٠٢٥"
$ws.Range("F46").Value = "١:٠٠:٥٧ م"

Write-Output "rows added"
